$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# N2 is stored as text (not a real date); Excel keeps string values as text automatically
$ws.Range("N2").Value = "2019-12-31 00:00:00"

$ws.Range("O2").Value = 2029905723
$ws.Range("P2").Value = 402279216.15
$ws.Range("Q2").Value = 424924472.3
$ws.Range("R2").Value = 1.3550770347
$ws.Range("S2").Value = 230242867.45
$ws.Range("T2").Value = -27.8480912565
$ws.Range("U2").Value = 646774199.84
$ws.Range("V2").Value = -11.2652041509
$ws.Range("W2").Value = 976483491.35
$ws.Range("X2").Value = 210679091.63
$ws.Range("Y2").Value = 41.1021219835
$ws.Range("Z2").Value = 82305628.89
$ws.Range("AA2").Value = -29.8179305584
$ws.Range("AB2").Value = 1053422231.65
$ws.Range("AC2").Value = 12.9649034039
$ws.Range("AD2").Value = -10.171094324
$ws.Range("AE2").Value = -26.4266801371
$ws.Range("AF2").Value = 145.6171519351
$ws.Range("AG2").Value = 48.104869122
